# The workbook gained one new weekly price record. In the canonical data
# this shows up as a new row inserted at row 43 (pushing every following
# row down by one, with the former last row, 105, becoming row 106).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 43; everything below shifts down by one.
$ws.Rows("43:43").Insert()

# Populate the new record.
$ws.Range("A43").Value = 6
$ws.Range("B43").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C43").Value = 'Metropolitana'
$ws.Range("D43").Value = 44482
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 100112001
$ws.Range("G43").Value = 'Berenjena'
$ws.Range("H43").Value = 'Sin especificar'
$ws.Range("I43").Value = 'Primera'
$ws.Range("J43").Value = 210
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 9000
$ws.Range("M43").Value = 8571
$ws.Range("N43").Value = '$/caja 60 unidades'
$ws.Range("O43").Value = 'Provincia de Huasco'
$ws.Range("P43").Value = 143
$ws.Range("Q43").Value = 60
$ws.Range("R43").Value = 'Hortaliza'
